# Manuscript revisions and response to reviewers; seedling plot analysis and
# figures in response to reviewers (round 2)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1 ("Seedling plot descriptions") content edits
# ---------------------------------------------------------------------

# Row 18/19 (Racetrack) canopy height + species text update
$ws1.Range("D18").Value = "10-18m, (Cycas micronesica 3m)"
$ws1.Range("D18").Characters(10, 18).Font.Italic = $true

$ws1.Range("C19").Value = "Ochrosia oppositifolia, Meiogyne cylindrocarpa, (Carica papaya) (Macaranga thompsonii, out of plot, but canopy overhanging plot)"
$ws1.Range("C19").Characters(1, 63).Font.Italic = $true
$ws1.Range("C19").Characters(66, 21).Font.Italic = $true

# Row 24 (Anao South) canopy height update
$ws1.Range("D24").Value = "9-15m, Cycas micronesica <3m"
$ws1.Range("D24").Characters(8, 17).Font.Italic = $true

# Row 34 (South Blas) adult tree count
$ws1.Range("D34").Value = 6

# Canopy cover percentage formatting fix (row 29, Racetrack)
$ws1.Range("C29").NumberFormat = "0%"
$ws1.Range("D29").NumberFormat = "0%"

# Update the current selection/cursor position on Sheet1
$ws1.Range("D6").Select()

# ---------------------------------------------------------------------
# Add Sheet2 with a small summary table header
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Site"
$ws2.Range("B1").Value = "adult trees"
$ws2.Range("C1").Value = "avg dbh"
$ws2.Range("D1").Value = "canopy cover"
$ws2.Range("E1").Value = "canopy height"

$ws2.Range("A2").Select()
$ws2.Activate()
